# Updates cryptos list values (Price / Volume(1h)) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.678.44'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').Value = '1.795.63'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''226.89'
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('E6').Value = '  +1.98%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '''32.92'
$ws.Range('E8').Value = '  +3.40%  '
$ws.Range('D9').Value = '''0.298'
$ws.Range('E9').Value = '  +2.18%  '
$ws.Range('D10').Value = '''0.0695'
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('D11').Value = '''0.0951'
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D12').Value = '2.055.14'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').Value = '''11.14'
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').Value = '1.801.55'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('E15').Value = '  +2.28%  '
$ws.Range('D16').Value = '34.580.97'
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').Value = '''4.29'
$ws.Range('E17').Value = '  +2.54%  '
$ws.Range('D18').Value = '''68.94'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('D19').Value = '''248.41'
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').Value = '0.0₃0803'
$ws.Range('E20').Value = '  +3.16%  '
$ws.Range('D21').Value = '''11.28'
$ws.Range('E21').Value = '  +2.97%  '
$ws.Range('D22').Value = '''1.00'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('E23').Value = '  +2.13%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '''2.07'
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''165.25'
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('D26').Value = '''7.27'
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('D27').Value = '''16.59'
$ws.Range('E27').Value = '  +1.67%  '
$ws.Range('E28').Value = '  +2.53%  '
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('D30').Value = '''4.15'
$ws.Range('E30').Value = '  +14.35%  '
$ws.Range('D31').Value = '''3.82'
$ws.Range('E31').Value = '  +3.41%  '
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('D33').Value = '''0.0524'
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('D34').Value = '''1.84'
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('D35').Value = '1.426.72'
$ws.Range('E35').Value = '  -1.35%  '
$ws.Range('D36').Value = '''2.58'
$ws.Range('E36').Value = '  +5.93%  '
$ws.Range('E37').Value = '  +3.33%  '
$ws.Range('E38').Value = '  +2.00%  '
$ws.Range('E39').Value = '  +0.65%  '
$ws.Range('D40').Value = '''85.66'
$ws.Range('E40').Value = '  +6.69%  '
$ws.Range('E41').Value = '  +0.82%  '
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('D43').Value = '''2.75'
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range('D44').Value = '''13.65'
$ws.Range('E44').Value = '  +1.10%  '
$ws.Range('D45').Value = '''0.0526'
$ws.Range('E45').Value = '  +3.64%  '
$ws.Range('D46').Value = '''6.11'
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('E47').Value = '  +0.67%  '
$ws.Range('D48').Value = '1.954.63'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('D49').Value = '''106.24'
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('E51').Value = '  -4.95%  '
